# Apply updated crypto price/volume figures (Sept 25 2024 GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.257.15'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('D3').Value = '2.580.07'
$ws.Range('E3').Value = '  -2.29%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''589.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.96%  '
$ws.Range('D6').Value = '''150.08'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.86%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('D10').Value = '''5.70'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.69%  '
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('E12').Value = '  -0.78%  '
$ws.Range('D13').Value = '''27.49'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.31%  '
$ws.Range('D14').Value = '3.044.22'
$ws.Range('E14').Value = '  -2.25%  '
$ws.Range('D15').Value = '63.073.01'
$ws.Range('E15').Value = '  -0.86%  '
$ws.Range('E16').Value = '  +5.17%  '
$ws.Range('D17').Value = '2.596.44'
$ws.Range('E17').Value = '  -1.83%  '
$ws.Range('D18').Value = '''12.24'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.23%  '
$ws.Range('E19').Value = '  +4.22%  '
$ws.Range('D20').Value = '''344.74'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.57%  '
$ws.Range('D21').Value = '''6.87'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.47%  '
$ws.Range('D22').Value = '''0.998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').Value = '''67.22'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('E24').Value = '  +1.15%  '
$ws.Range('D25').Value = '''9.26'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.36%  '
$ws.Range('D26').Value = '''1.67'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.09%  '
$ws.Range('D27').Value = '''564.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.47%  '
$ws.Range('D28').Value = '''8.05'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.84%  '
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('E31').Value = '  -1.22%  '
$ws.Range('D32').Value = '0.0₃0845'
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('D33').Value = '''1.76'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D34').Value = '''5.22'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.40%  '
$ws.Range('D35').Value = '''166.71'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.55%  '
$ws.Range('D36').Value = '''0.412'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.66%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').Value = '''19.45'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.53%  '
$ws.Range('E39').Value = '  -0.46%  '
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').Value = '''166.82'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.09%  '
$ws.Range('D42').Value = '''39.51'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.44%  '
$ws.Range('D43').Value = '''3.95'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.98%  '
$ws.Range('D44').Value = '''22.71'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.79%  '
$ws.Range('E45').Value = '  +2.90%  '
$ws.Range('D46').Value = '''2.07'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.98%  '
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('D48').Value = '''0.0252'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.99%  '
$ws.Range('E49').Value = '  +0.53%  '
$ws.Range('D50').Value = '''19.03'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.19%  '
$ws.Range('D51').Value = '0.0₆0232'
$ws.Range('E51').Value = '  +17.60%  '
